$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5181.222
$ws.Range("I43").Value = 4168.3335
$ws.Range("J43").Value = 6194.1113
$ws.Range("K43").Value = 4168.3335
$ws.Range("L43").Value = 6194.1113
$ws.Range("M43").Value = -4099.3335
$ws.Range("N43").Value = -6332.1113

$ws.Range("H99").Value = 421.42856
$ws.Range("I99").Value = 375.25
$ws.Range("J99").Value = 698.5
$ws.Range("K99").Value = 1125.75
$ws.Range("L99").Value = 2095.5
$ws.Range("M99").Value = 372.25
$ws.Range("N99").Value = -5091.5

$ws.Range("H116").Value = 8175.625
$ws.Range("I116").Value = 7681
$ws.Range("J116").Value = 9000
$ws.Range("K116").Value = 7681
$ws.Range("L116").Value = 9000
$ws.Range("M116").Value = -4239
$ws.Range("N116").Value = -15884

$ws.Range("H137").Value = 3298.2742
$ws.Range("J137").Value = 5720.1113
$ws.Range("L137").Value = 17160.3339
$ws.Range("N137").Value = -22260.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 688.4545000000001
$ws.Range("I2").Value = 569.2778
$ws.Range("K2").Value = 569.2778
$ws.Range("M2").Value = -456.2778

$ws.Range("H32").Value = 11116801
$ws.Range("I32").Value = 14708419
$ws.Range("K32").Value = 14708419
$ws.Range("M32").Value = -14708132

$ws.Range("H45").Value = 31252846
$ws.Range("I45").Value = 55557372
$ws.Range("K45").Value = 55557372
$ws.Range("M45").Value = -55556995

$ws.Range("H116").Value = 688.4545000000001
$ws.Range("I116").Value = 569.2778
$ws.Range("K116").Value = 569.2778
$ws.Range("M116").Value = 1724.7222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 688.4545000000001
$ws.Range("I3").Value = 569.2778
$ws.Range("K3").Value = 569.2778
$ws.Range("M3").Value = -455.2778

$ws.Range("H20").Value = 4208.25
$ws.Range("I20").Value = 4597.875
$ws.Range("K20").Value = 4597.875
$ws.Range("M20").Value = -4350.875

$ws.Range("H99").Value = 2548.48
$ws.Range("I99").Value = 1845.2
$ws.Range("J99").Value = 3603.4
$ws.Range("K99").Value = 1845.2
$ws.Range("L99").Value = 3603.4
$ws.Range("M99").Value = -347.2
$ws.Range("N99").Value = -6599.4

$ws.Range("H105").Value = 3025.7222
$ws.Range("I105").Value = 4243.6665
$ws.Range("K105").Value = 4243.6665
$ws.Range("M105").Value = -2496.6665

$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 54056.5
$ws.Range("J124").Value = 54056.5
$ws.Range("L124").Value = 54056.5
$ws.Range("N124").Value = -58966.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 423.2
$ws.Range("I86").Value = 396.5
$ws.Range("J86").Value = 463.25
$ws.Range("K86").Value = 1189.5
$ws.Range("L86").Value = 1389.75
$ws.Range("M86").Value = -3.5
$ws.Range("N86").Value = -3761.75

$ws.Range("H89").Value = 423.2
$ws.Range("I89").Value = 396.5
$ws.Range("J89").Value = 463.25
$ws.Range("K89").Value = 3568.5
$ws.Range("L89").Value = 4169.25
$ws.Range("M89").Value = 2359.5
$ws.Range("N89").Value = -16025.25

$ws.Range("H131").Value = 11542.533
$ws.Range("J131").Value = 11542.533
$ws.Range("L131").Value = 34627.599
$ws.Range("N131").Value = -44707.599

$ws.Range("H138").Value = 4307.615
$ws.Range("I138").Value = 2749.75
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 8249.25
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = -3109.25
$ws.Range("N138").Value = -25280

$ws.Range("H140").Value = 3080.2778
$ws.Range("I140").Value = 3090.3125
$ws.Range("K140").Value = 9270.9375
$ws.Range("M140").Value = -4090.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 29999.75
$ws.Range("I70").Value = 36666.332
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 36666.332
$ws.Range("L70").Value = 10000
$ws.Range("M70").Value = -36396.332
$ws.Range("N70").Value = -10540

$ws.Range("H73").Value = 29999.75
$ws.Range("I73").Value = 36666.332
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 36666.332
$ws.Range("L73").Value = 10000
$ws.Range("M73").Value = -35730.332
$ws.Range("N73").Value = -11872

$ws.Range("H107").Value = 761.38464
$ws.Range("I107").Value = 747.4286
$ws.Range("J107").Value = 777.6667
$ws.Range("K107").Value = 747.4286
$ws.Range("L107").Value = 777.6667
$ws.Range("M107").Value = 1172.5714
$ws.Range("N107").Value = -4617.6667

$ws.Range("H122").Value = 1307.6818
$ws.Range("I122").Value = 1352.7778
$ws.Range("J122").Value = 1104.75
$ws.Range("K122").Value = 4058.3334
$ws.Range("L122").Value = 3314.25
$ws.Range("M122").Value = -1608.3334
$ws.Range("N122").Value = -8214.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3046.1765
$ws.Range("I46").Value = 2291
$ws.Range("J46").Value = 4430.6665
$ws.Range("K46").Value = 2291
$ws.Range("L46").Value = 4430.6665
$ws.Range("M46").Value = -2103
$ws.Range("N46").Value = -4806.6665

$ws.Range("H50").Value = 40000
$ws.Range("J50").Value = 40000
$ws.Range("L50").Value = 40000
$ws.Range("N50").Value = -41274

$ws.Range("H93").Value = 333334660
$ws.Range("I93").Value = 500001000
$ws.Range("K93").Value = 500001000
$ws.Range("M93").Value = -499999752

$ws.Range("H100").Value = 4361.75
$ws.Range("I100").Value = 3975
$ws.Range("K100").Value = 3975
$ws.Range("M100").Value = -3434

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8262.75
$ws.Range("I81").Value = 2240.2
$ws.Range("J81").Value = 18300.334
$ws.Range("K81").Value = 4480.4
$ws.Range("L81").Value = 36600.668
$ws.Range("M81").Value = -3419.4
$ws.Range("N81").Value = -38722.668

$ws.Range("H84").Value = 8262.75
$ws.Range("I84").Value = 2240.2
$ws.Range("J84").Value = 18300.334
$ws.Range("K84").Value = 22402
$ws.Range("L84").Value = 183003.34
$ws.Range("M84").Value = -17098
$ws.Range("N84").Value = -193611.34

$ws.Range("H100").Value = 573.37036
$ws.Range("I100").Value = 401.6
$ws.Range("J100").Value = 674.41174
$ws.Range("K100").Value = 803.2
$ws.Range("L100").Value = 1348.82348
$ws.Range("M100").Value = -262.2
$ws.Range("N100").Value = -2430.82348

$ws.Range("H107").Value = 12821410
$ws.Range("I107").Value = 17242380
$ws.Range("K107").Value = 51727140
$ws.Range("M107").Value = -51725220

$ws.Range("H119").Value = 112750
$ws.Range("I119").Value = 67500
$ws.Range("J119").Value = 158000
$ws.Range("K119").Value = 67500
$ws.Range("L119").Value = 158000
$ws.Range("M119").Value = -62662
$ws.Range("N119").Value = -167676

$ws.Range("H122").Value = 2367.7917
$ws.Range("I122").Value = 2196.1
$ws.Range("J122").Value = 3226.25
$ws.Range("K122").Value = 6588.299999999999
$ws.Range("L122").Value = 9678.75
$ws.Range("M122").Value = -4138.299999999999
$ws.Range("N122").Value = -14578.75

$ws.Range("H126").Value = 1426.5714
$ws.Range("I126").Value = 1426.5714
$ws.Range("K126").Value = 4279.7142
$ws.Range("M126").Value = -1809.7142

$ws.Range("H136").Value = 3946.2222
$ws.Range("I136").Value = 3883
$ws.Range("K136").Value = 11649
$ws.Range("M136").Value = -9099
